# Rename the "type" column in the Variables table to "variable-type".
# Renaming the header cell's value also updates the bound ListObject's
# table column name (xl/tables/table3.xml) and causes the shared string
# "type" to be dropped / "variable-type" to be (re)created in
# xl/sharedStrings.xml.
$wb = $excel.ActiveWorkbook
$wsVar = $wb.Worksheets.Item("Variables")
$wsVar.Range("D1").Value = "variable-type"

# Reflect the author's final selection/active-sheet state: the active
# sheet moves from "Cells" to "Variables", with cell E8 selected there.
$wsVar.Activate()
$wsVar.Range("E8").Select()
